$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44477
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 6000
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = 6000
$ws.Range("P2").Value = 375
$ws.Range("D3").Value = 44403
$ws.Range("J3").Value = 43
$ws.Range("K3").Value = 6000
$ws.Range("L3").Value = 6000
$ws.Range("M3").Value = 6000
$ws.Range("P3").Value = 375
$ws.Range("D5").Value = 44438
$ws.Range("J5").Value = 34
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = 5500
$ws.Range("P5").Value = 344
$ws.Range("D6").Value = 44363
$ws.Range("J6").Value = 160
$ws.Range("K6").Value = 5500
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 5750
$ws.Range("P6").Value = 359
$ws.Range("D7").Value = 44575
$ws.Range("J7").Value = 61
$ws.Range("K7").Value = 8000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 8000
$ws.Range("P7").Value = 500
$ws.Range("D8").Value = 44474
$ws.Range("J8").Value = 52
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = 5500
$ws.Range("P8").Value = 344
$ws.Range("D9").Value = 44355
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 6000
$ws.Range("L9").Value = 6000
$ws.Range("M9").Value = 6000
$ws.Range("P9").Value = 375
$ws.Range("D10").Value = 44371
$ws.Range("J10").Value = 34
$ws.Range("K10").Value = 5500
$ws.Range("L10").Value = 6000
$ws.Range("M10").Value = 5750
$ws.Range("P10").Value = 359
$ws.Range("D11").Value = 44582
$ws.Range("J11").Value = 52
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("P11").Value = 438
$ws.Range("D12").Value = 44698
$ws.Range("J12").Value = 34
$ws.Range("K12").Value = 6000
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 6500
$ws.Range("P12").Value = 406
$ws.Range("D13").Value = 44376
$ws.Range("J13").Value = 43
$ws.Range("K13").Value = 4500
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = 4756
$ws.Range("P13").Value = 297
$ws.Range("D14").Value = 44341
$ws.Range("J14").Value = 51
$ws.Range("K14").Value = 5500
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = 5755
$ws.Range("P14").Value = 360
$ws.Range("D15").Value = 44573
$ws.Range("J15").Value = 34
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 8000
$ws.Range("M15").Value = 8000
$ws.Range("P15").Value = 500
$ws.Range("D16").Value = 44715
$ws.Range("J16").Value = 70
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 6000
$ws.Range("M16").Value = 5500
$ws.Range("P16").Value = 344
$ws.Range("D17").Value = 44455
$ws.Range("J17").Value = 52
$ws.Range("K17").Value = 5000
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = 5500
$ws.Range("P17").Value = 344
$ws.Range("D18").Value = 44308
$ws.Range("J18").Value = 70
$ws.Range("K18").Value = 6000
$ws.Range("L18").Value = 6000
$ws.Range("M18").Value = 6000
$ws.Range("P18").Value = 375
$ws.Range("D19").Value = 44306
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 6000
$ws.Range("L19").Value = 6000
$ws.Range("M19").Value = 6000
$ws.Range("P19").Value = 375
$ws.Range("D20").Value = 44330
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = 6000
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = 6000
$ws.Range("P20").Value = 375
$ws.Range("D21").Value = 44407
$ws.Range("J21").Value = 45
$ws.Range("K21").Value = 5500
$ws.Range("L21").Value = 6000
$ws.Range("M21").Value = 5744
$ws.Range("P21").Value = 359
$ws.Range("D22").Value = 44589
$ws.Range("J22").Value = 52
$ws.Range("K22").Value = 8000
$ws.Range("L22").Value = 8000
$ws.Range("M22").Value = 8000
$ws.Range("P22").Value = 500
$ws.Range("D23").Value = 44691
$ws.Range("J23").Value = 61
$ws.Range("K23").Value = 6000
$ws.Range("L23").Value = 7000
$ws.Range("M23").Value = 6508
$ws.Range("P23").Value = 407
$ws.Range("D24").Value = 44467
$ws.Range("J24").Value = 52
$ws.Range("K24").Value = 5000
$ws.Range("L24").Value = 6000
$ws.Range("M24").Value = 5500
$ws.Range("P24").Value = 344
$ws.Range("D25").Value = 44358
$ws.Range("J25").Value = 52
$ws.Range("K25").Value = 6000
$ws.Range("L25").Value = 6000
$ws.Range("M25").Value = 6000
$ws.Range("P25").Value = 375
$ws.Range("D26").Value = 44313
$ws.Range("J26").Value = 34
$ws.Range("K26").Value = 6000
$ws.Range("L26").Value = 6000
$ws.Range("M26").Value = 6000
$ws.Range("P26").Value = 375
$ws.Range("D27").Value = 44350
$ws.Range("J27").Value = 25
$ws.Range("K27").Value = 6000
$ws.Range("L27").Value = 6000
$ws.Range("M27").Value = 6000
$ws.Range("P27").Value = 375
$ws.Range("D28").Value = 44328
$ws.Range("J28").Value = 160
$ws.Range("K28").Value = 6000
$ws.Range("L28").Value = 6000
$ws.Range("M28").Value = 6000
$ws.Range("P28").Value = 375
